$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2 with new user data
$ws.Range("B2").Value = "aa"
$ws.Range("C2").Value = "aa@gmail.com"
$ws.Range("D2").Value = '$2b$10$aMWseoSQNXof5F78tCDz6uN.xy/U1H3cvPRaPbjGKw0/NsjlvE55O'
$ws.Range("E2").Value = "/public/images/4.jpg"
$ws.Range("F2").Value = 1397
$ws.Range("G2").Value = $false

# Remove row 3 entirely (second user record)
$ws.Range("A3:G3").EntireRow.Delete()
